$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
$ws.Range("H70").Value = 1803.4546
$ws.Range("J70").Value = 1772.2222
$ws.Range("L70").Value = 5316.6666
$ws.Range("N70").Value = -5856.6666
$ws.Range("H73").Value = 1803.4546
$ws.Range("J73").Value = 1772.2222
$ws.Range("L73").Value = 5316.6666
$ws.Range("N73").Value = -7188.6666
$ws.Range("H88").Value = 9477.799999999999
$ws.Range("J88").Value = 9872.25
$ws.Range("L88").Value = 9872.25
$ws.Range("N88").Value = -10684.25
$ws.Range("H91").Value = 9477.799999999999
$ws.Range("J91").Value = 9872.25
$ws.Range("L91").Value = 9872.25
$ws.Range("N91").Value = -12680.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3033.889
$ws.Range("I88").Value = 2099.8
$ws.Range("J88").Value = 4201.5
$ws.Range("K88").Value = 2099.8
$ws.Range("L88").Value = 4201.5
$ws.Range("M88").Value = -1693.8
$ws.Range("N88").Value = -5013.5
$ws.Range("H91").Value = 3033.889
$ws.Range("I91").Value = 2099.8
$ws.Range("J91").Value = 4201.5
$ws.Range("K91").Value = 2099.8
$ws.Range("L91").Value = 4201.5
$ws.Range("M91").Value = -695.8000000000002
$ws.Range("N91").Value = -7009.5
$ws.Range("H102").Value = 2883.6428
$ws.Range("I102").Value = 2861.7273
$ws.Range("K102").Value = 2861.7273
$ws.Range("M102").Value = -1239.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 12071.833
$ws.Range("I20").Value = 2204
$ws.Range("J20").Value = 17005.75
$ws.Range("K20").Value = 2204
$ws.Range("L20").Value = 17005.75
$ws.Range("M20").Value = -1957
$ws.Range("N20").Value = -17499.75
$ws.Range("H86").Value = 18833.334
$ws.Range("I86").Value = 1500
$ws.Range("K86").Value = 1500
$ws.Range("M86").Value = -377
$ws.Range("H89").Value = 18833.334
$ws.Range("I89").Value = 1500
$ws.Range("K89").Value = 7500
$ws.Range("M89").Value = -1884
$ws.Range("H134").Value = 3415.5454
$ws.Range("J134").Value = 1300
$ws.Range("L134").Value = 3900
$ws.Range("N134").Value = -8970

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 19926.5
$ws.Range("J62").Value = 25002
$ws.Range("L62").Value = 25002
$ws.Range("N62").Value = -26250
$ws.Range("H65").Value = 19926.5
$ws.Range("J65").Value = 25002
$ws.Range("L65").Value = 125010
$ws.Range("N65").Value = -131250
$ws.Range("H99").Value = 3012.3333
$ws.Range("I99").Value = 2814.8
$ws.Range("K99").Value = 2814.8
$ws.Range("M99").Value = -1316.8
$ws.Range("H107").Value = 1770.7142
$ws.Range("I107").Value = 1799.3334
$ws.Range("J107").Value = 1749.25
$ws.Range("K107").Value = 1799.3334
$ws.Range("L107").Value = 1749.25
$ws.Range("M107").Value = 120.6666
$ws.Range("N107").Value = -5589.25
$ws.Range("H126").Value = 3012.3333
$ws.Range("I126").Value = 2814.8
$ws.Range("K126").Value = 8444.400000000001
$ws.Range("M126").Value = -5974.400000000001
$ws.Range("H132").Value = 4255.3335
$ws.Range("I132").Value = 3775
$ws.Range("K132").Value = 11325
$ws.Range("M132").Value = -8795

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 194.5
$ws.Range("I38").Value = 240
$ws.Range("K38").Value = 720
$ws.Range("M38").Value = -373
$ws.Range("H92").Value = 297.25
$ws.Range("I92").Value = 362
$ws.Range("J92").Value = 103
$ws.Range("K92").Value = 1086
$ws.Range("L92").Value = 309
$ws.Range("M92").Value = 162
$ws.Range("N92").Value = -2805
$ws.Range("H109").Value = 750
$ws.Range("I109").Value = 750
$ws.Range("K109").Value = 2250
$ws.Range("M109").Value = -1210
$ws.Range("H128").Value = 299999
$ws.Range("I128").Value = 299999
$ws.Range("K128").Value = 899997
$ws.Range("M128").Value = -895017
$ws.Range("H139").Value = 112788.89
$ws.Range("I139").Value = 112788.89
$ws.Range("K139").Value = 338366.67
$ws.Range("M139").Value = -333226.67

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10857
$ws.Range("I80").Value = 3666.3333
$ws.Range("J80").Value = 16250
$ws.Range("K80").Value = 3666.3333
$ws.Range("L80").Value = 16250
$ws.Range("M80").Value = -2668.3333
$ws.Range("N80").Value = -18246
$ws.Range("H83").Value = 10857
$ws.Range("I83").Value = 3666.3333
$ws.Range("J83").Value = 16250
$ws.Range("K83").Value = 18331.6665
$ws.Range("L83").Value = 81250
$ws.Range("M83").Value = -13339.6665
$ws.Range("N83").Value = -91234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2977.1428
$ws.Range("I68").Value = 2977.1428
$ws.Range("K68").Value = 2977.1428
$ws.Range("M68").Value = -2228.1428
$ws.Range("H71").Value = 2977.1428
$ws.Range("I71").Value = 2977.1428
$ws.Range("K71").Value = 14885.714
$ws.Range("M71").Value = -11141.714
$ws.Range("H76").Value = 24234.25
$ws.Range("J76").Value = 24234.25
$ws.Range("L76").Value = 24234.25
$ws.Range("N76").Value = -24910.25
$ws.Range("H79").Value = 24234.25
$ws.Range("J79").Value = 24234.25
$ws.Range("L79").Value = 24234.25
$ws.Range("N79").Value = -26574.25
$ws.Range("H94").Value = 39987.5
$ws.Range("J94").Value = 39987.5
$ws.Range("L94").Value = 39987.5
$ws.Range("N94").Value = -41339.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4946.7144
$ws.Range("I81").Value = 1725.4
$ws.Range("J81").Value = 13000
$ws.Range("K81").Value = 3450.8
$ws.Range("L81").Value = 26000
$ws.Range("M81").Value = -2389.8
$ws.Range("N81").Value = -28122
$ws.Range("H84").Value = 4946.7144
$ws.Range("I84").Value = 1725.4
$ws.Range("J84").Value = 13000
$ws.Range("K84").Value = 17254
$ws.Range("L84").Value = 130000
$ws.Range("M84").Value = -11950
$ws.Range("N84").Value = -140608
